# Modulos_Portal-Asignaciones - v001.xlsx
# Se crean archivos con funciones genericas para el consumo del webservice por modulo.
#
# - Adds a new shared string "Enrique" and uses it to populate the "BD"
#   column (H) for the "Normal" (row 20) and "Corte" (row 21) manifest rows,
#   matching the existing "Documentar NUI" row (14) which already carries an
#   analogous "Back" (I) assignment of "Gustavo".
# - Row 14's "BD" cell (H14) gets the new "Enrique" assignment too.
# - Rows 20/21 also get the existing "Gustavo" shared string placed in the
#   "Back" column (I), mirroring row 14's layout.
# - The active selection moves from M9 to F16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content ------------------------------------------------------

# H14 was an empty, styled placeholder cell -> now holds "Enrique".
$ws.Range("H14").Value = "Enrique"

# H20/I20 and H21/I21 were empty placeholder cells (italic style, no
# content) -> now hold "Enrique" / "Gustavo" respectively, matching the
# centered, non-italic style already used by H14/I14.
$ws.Range("H20").Value = "Enrique"
$ws.Range("I20").Value = "Gustavo"

$ws.Range("H21").Value = "Enrique"
$ws.Range("I21").Value = "Gustavo"

# The target style (already used by H14/I14) is identical to the cells'
# current style except that it is not italic, so flipping Italic off
# reuses that existing style rather than creating a new one.
$ws.Range("H20:I21").Font.Italic = $false

# --- View state ----------------------------------------------------------
# Selection moves from M9 to F16; also nudge the scroll position toward
# the new top-left cell (E8) for best effort parity with the source view.
$ws.Range("F16").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 5
